$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback identifiers / hashes / timestamps replacing the old ones.
# ---------------------------------------------------------------------------
$oldId1 = "0ad5e829-c9a0-43b0-8bab-1694ebf7776b"
$newId1 = "123824c4-9ff1-495b-9e17-1bca96978028"
$oldId2 = "36c11795-e69f-455c-82d3-d3b6aa28a8ea"
$newId2 = "ffff8411bb41-5ec9-4414-82c8-0b4818238600"
$newHash = "2045dbcd49a61fdcc75759731960b0993c9f0fb1"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = "2016-08-30 05:02:41"

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = "2016-08-30 05:02:41"

# Hyperlinks keep pointing at the same targets, only the displayed text
# (the file name) changes, so the links are rebuilt in place.
$hlTargetB2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId1.md"
$hlTargetB3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId2.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlTargetB2, [Type]::Missing, [Type]::Missing, "e2e\$newId1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hlTargetB3, [Type]::Missing, [Type]::Missing, "e2e\$newId2.md")

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("G2").Value = "$newId1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 05:02:36"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("J2").Value = "$newId1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 05:02:53"

$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("G3").Value = "$newId1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 05:02:36"
$wsZhCn.Range("I3").Value = "$newId2.md"
$wsZhCn.Range("J3").Value = "$newId1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 05:02:53"

$hlTargetA2zh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId1.md"
$hlTargetI2zh = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ab6813579a0379197b1755c7ca0fdedc80583866/e2e/$oldId1.md"
$hlTargetA3zh = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId2.md"
$hlTargetI3zh = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ab6813579a0379197b1755c7ca0fdedc80583866/e2e/$oldId2.md"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hlTargetA2zh, [Type]::Missing, [Type]::Missing, "$newId1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $hlTargetI2zh, [Type]::Missing, [Type]::Missing, "$newId1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hlTargetA3zh, [Type]::Missing, [Type]::Missing, "$newId2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $hlTargetI3zh, [Type]::Missing, [Type]::Missing, "$newId2.md")

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("G2").Value = "$newId1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 05:02:41"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("J2").Value = "$newId1.$newHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 05:03:00"

$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("G3").Value = "$newId1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 05:02:41"
$wsDeDe.Range("I3").Value = "$newId2.md"
$wsDeDe.Range("J3").Value = "$newId1.$newHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 05:03:00"

$hlTargetA2de = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId1.md"
$hlTargetI2de = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1c8bf4d1082bac4eb03911e9e72cdc4420da1e44/e2e/$oldId1.md"
$hlTargetA3de = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/87e7f838f6e14b813d7ec7f120afa372f5ac8d6b/e2e/$oldId2.md"
$hlTargetI3de = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1c8bf4d1082bac4eb03911e9e72cdc4420da1e44/e2e/$oldId2.md"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hlTargetA2de, [Type]::Missing, [Type]::Missing, "$newId1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $hlTargetI2de, [Type]::Missing, [Type]::Missing, "$newId1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hlTargetA3de, [Type]::Missing, [Type]::Missing, "$newId2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $hlTargetI3de, [Type]::Missing, [Type]::Missing, "$newId2.md")
